# Update "想去人数" (interested-count) values in column F
# for the "展览" and "全部类型" worksheets, reflecting newly scraped totals.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 59
$wsExhibit.Range("F3").Value  = 787
$wsExhibit.Range("F5").Value  = 58
$wsExhibit.Range("F6").Value  = 73
$wsExhibit.Range("F7").Value  = 273
$wsExhibit.Range("F8").Value  = 3933
$wsExhibit.Range("F10").Value = 4620
$wsExhibit.Range("F11").Value = 509
$wsExhibit.Range("F12").Value = 1166
$wsExhibit.Range("F13").Value = 74

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 59
$wsAll.Range("F3").Value  = 787
$wsAll.Range("F5").Value  = 58
$wsAll.Range("F6").Value  = 73
$wsAll.Range("F8").Value  = 273
$wsAll.Range("F9").Value  = 3933
$wsAll.Range("F11").Value = 4620
$wsAll.Range("F12").Value = 509
$wsAll.Range("F13").Value = 1166
$wsAll.Range("F14").Value = 74
